$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.849948
$ws.Range("H2").Value = 599.00257804
$ws.Range("M2").Value = 0.7396380624240296
$ws.Range("N2").Value = 71.24459305997657
$ws.Range("G3").Value = 11.82648092
$ws.Range("H3").Value = 1072.18224382
$ws.Range("M3").Value = 1.211473474975781
$ws.Range("N3").Value = 198.8988328486925
$ws.Range("G4").Value = 3.77786168
$ws.Range("H4").Value = 117.8756196
$ws.Range("M4").Value = 0.5679745803828722
$ws.Range("N4").Value = 29.91926874613998
$ws.Range("G5").Value = 4.24773734
$ws.Range("H5").Value = 208.2055144
$ws.Range("M5").Value = 0.5454160940460165
$ws.Range("N5").Value = 51.6332671987437
$ws.Range("G6").Value = 1.09936002
$ws.Range("H6").Value = 17.09028364
$ws.Range("M6").Value = 0.2682831429041576
$ws.Range("N6").Value = 7.407689108259147
$ws.Range("G7").Value = 1.52148366
$ws.Range("H7").Value = 39.70372676
$ws.Range("M7").Value = 0.3007313816718317
$ws.Range("N7").Value = 14.59971592972864
$ws.Range("G8").Value = 0.56847474
$ws.Range("H8").Value = 5.76777984
$ws.Range("M8").Value = 0.2045249537703454
$ws.Range("N8").Value = 3.269834813216015
$ws.Range("G9").Value = 0.71084504
$ws.Range("H9").Value = 12.56477354
$ws.Range("M9").Value = 0.1726946855599201
$ws.Range("N9").Value = 5.778872800209819
$ws.Range("G10").Value = 0.28976586
$ws.Range("H10").Value = 2.233636
$ws.Range("M10").Value = 0.1197891019122118
$ws.Range("N10").Value = 1.459462870458888
$ws.Range("G11").Value = 0.3943485
$ws.Range("H11").Value = 5.443070659999999
$ws.Range("M11").Value = 0.113917157464007
$ws.Range("N11").Value = 3.123867325064587
$ws.Range("G12").Value = 0.1724858
$ws.Range("H12").Value = 1.09011984
$ws.Range("M12").Value = 0.07780214060117015
$ws.Range("N12").Value = 0.7553301821596546
$ws.Range("G13").Value = 0.24672746
$ws.Range("H13").Value = 2.92827752
$ws.Range("M13").Value = 0.08703414882603347
$ws.Range("N13").Value = 2.351381151570151

